# Suivi LDLC smartphones: append a new price-history snapshot column.
#
# The sheet has one header row (row 1) with timestamped snapshot columns
# running up to CZ, followed by two fixed columns "nom" (product name) and
# "url_produit" (product URL). This edit inserts one new snapshot column
# right before "nom"/"url_produit":
#   - DA1 gets the new snapshot timestamp header.
#   - "nom"/"url_produit" (and all per-row data under them) shift one
#     column to the right, DA -> DB, DB -> DC.
#   - For rows that already had a price series (rows 2-80, where the last
#     snapshot column CZ holds a numeric price), the new DA column is
#     populated with that same latest price (i.e. the price carried
#     forward into the new snapshot).
#   - For rows with no price series yet (rows 81-206, CZ blank), the new
#     DA column is left blank too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newHeader = "2026-02-01 12:22:45"

# Column numbers: CZ = 104, DA = 105, DB = 106, DC = 107
$lastSnapshotCol = 104
$newCol = 105

# Find the last used row on the sheet (206 in this workbook).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Insert a new blank column before the current DA ("nom"), shifting
# "nom"/"url_produit" (and everything to their right) one column over.
$ws.Columns.Item($newCol).Insert()

# Header for the freshly inserted snapshot column.
$ws.Cells.Item(1, $newCol).Value = $newHeader

# Carry the latest known price (column CZ) forward into the new column
# for every data row that already has one.
for ($r = 2; $r -le $lastRow; $r++) {
    $price = $ws.Cells.Item($r, $lastSnapshotCol).Value()
    if ($price -ne "" -and $price -ne $null) {
        $ws.Cells.Item($r, $newCol).Value = $price
    }
}
